$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Common" folder rows (2-10): fill in Reviewer (C) / Status (D) ---
# These cells were previously blank (generic style). The Reviewer column's
# font needs to match the rest of the "Reviewer" column (Microsoft YaHei UI),
# same as the existing data further down the sheet, so it collapses onto the
# same cell style already used for B2:B10 instead of minting a new one.
$ws.Range("C2:C10").Value = "Baotong"
$ws.Range("C2:C10").Font.Name = "Microsoft YaHei UI"

$ws.Range("D2").Value = "Done"
$ws.Range("D3").Value = "private"
$ws.Range("D4").Value = "private"
$ws.Range("D5").Value = "Done"
$ws.Range("D6").Value = "Done"
$ws.Range("D7").Value = "Done"
$ws.Range("D8").Value = "Done"
$ws.Range("D9").Value = "Done"
$ws.Range("D10").Value = "Done"

# --- Row 18 (ParameterModel.js): reviewer changed from Jon to Baotong, ---
# --- and the status is now filled in as "private" ---
$ws.Range("C18").Value = "Baotong"
$ws.Range("D18").Value = "private"

# --- Row 37 (ReportViewerInitializer.js): status updated ---
$ws.Range("D37").Value = "Private"

# --- Row 38 (SubscriptionModel.js): status updated ---
$ws.Range("D38").Value = "Private"

# --- Restore the view: scroll back to the top and select D5 ---
$ws.Range("D5").Select()
